$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: set new value
$ws.Range("F2").Value = "33,33 TL - 33,33 TL"

# J6, K6: clear existing values
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""

# F7: set new value
$ws.Range("F7").Value = "%3"

# K12: clear existing value
$ws.Range("K12").Value = ""

# E13: update value
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"

# J13, K13: clear existing values
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""

# J14, K14: clear existing values
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
